$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date value in C2
$ws.Range("C2").Value = 45186

# Add friendly text as second argument to the HYPERLINK formulas
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/artfynd/A 14042-2023.xlsx", "A 14042-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/kartor/A 14042-2023.png", "A 14042-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/klagomål/A 14042-2023.docx", "A 14042-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/klagomålsmail/A 14042-2023.docx", "A 14042-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/tillsyn/A 14042-2023.docx", "A 14042-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PARTILLE/tillsynsmail/A 14042-2023.docx", "A 14042-2023")'
